# Applies the cryptocurrency price/volume updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values that look like plain decimal numbers must be forced back to
# Text format before assignment, otherwise Excel auto-converts them to numbers
# and the literal formatting (trailing zeros, thousand separators) is lost.
$textForceCells = @('D5', 'D11', 'D16', 'D20', 'D25', 'D27', 'D35', 'D37', 'D40', 'D43', 'D44', 'D46', 'D51')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.618.00'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '1.596.04'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '210.93'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range('D11').Value = '0.0836'
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').Value = '1.820.09'
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').Value = '1.611.39'
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').Value = '65.00'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = '26.605.69'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').Value = '0.0₃0737'
$ws.Range('E18').Value = '  +0.70%  '
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('D20').Value = '208.17'
$ws.Range('E20').Value = '  -0.86%  '
$ws.Range('E21').Value = '  +4.87%  '
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('D25').Value = '144.97'
$ws.Range('E25').Value = '  -1.20%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').Value = '7.12'
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('E29').Value = '  -0.54%  '
$ws.Range('E30').Value = '  +0.60%  '
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('E32').Value = '  -0.68%  '
$ws.Range('E33').Value = '  +0.40%  '
$ws.Range('D34').Value = '1.277.46'
$ws.Range('E34').Value = '  -1.26%  '
$ws.Range('D35').Value = '0.615'
$ws.Range('E35').Value = '  -8.16%  '
$ws.Range('E36').Value = '  +0.85%  '
$ws.Range('D37').Value = '1.49'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  -1.05%  '
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('D40').Value = '1.01'
$ws.Range('E40').Value = '  +16.94%  '
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('E42').Value = '  +0.20%  '
$ws.Range('D43').Value = '0.784'
$ws.Range('E43').Value = '  -1.00%  '
$ws.Range('D44').Value = '64.10'
$ws.Range('E44').Value = '  +0.37%  '
$ws.Range('D45').Value = '1.732.37'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = '90.07'
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('E47').Value = '  -1.50%  '
$ws.Range('E48').Value = '  +3.50%  '
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('D51').Value = '7.43'
$ws.Range('E51').Value = '  -1.21%  '
